# Updated symbol list on Sat Dec 17 23:32:59 UTC 2022 with GitHub Actions
#
# Applies refreshed coin prices plus the re-ordering of the coin rows 10-18
# (One/WazirX/MandalaExchangeToken/.../CoinExToken rotate one slot, with new
# quotes) and the BKEXToken/CEJI swap in rows 42-43.
#
# Price-like values are written with a leading apostrophe so Excel stores
# them as literal text (preserving exact digits/trailing zeros) instead of
# converting them to floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- simple price refreshes (rows 2-8) ---
$ws.Range("D2").Value  = "'239.94"
$ws.Range("D3").Value  = "'21.82"
$ws.Range("D4").Value  = "'5.411"
$ws.Range("D5").Value  = "'0.05574"
$ws.Range("D6").Value  = "'6.459"
$ws.Range("D7").Value  = "'3.360"
$ws.Range("D8").Value  = "'0.8047"

# --- row 10: was One -> now WazirX ---
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1413"
$ws.Range("E10").Value = "9WazirXWRX"

# --- row 11: was WazirX -> now MandalaExchangeToken ---
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07314"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

# --- row 12: was MandalaExchangeToken -> now LiechtensteinCryptoassetsExchange ---
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03284"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

# --- row 13: was LiechtensteinCryptoassetsExchange -> now BitrueCoin ---
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.02960"
$ws.Range("E13").Value = "12BitrueCoinBTR"

# --- row 14: was BitrueCoin -> now BitMartToken ---
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09251"
$ws.Range("E14").Value = "13BitMartTokenBMX"

# --- row 15: was BitMartToken -> now BitForexToken ---
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001667"
$ws.Range("E15").Value = "14BitForexTokenBF"

# --- row 16: was BitForexToken -> now MCDex ---
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.245"
$ws.Range("E16").Value = "15MCDexMCB"

# --- row 17: was MCDex -> now CoinExToken ---
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04760"
$ws.Range("E17").Value = "16CoinExTokenCET"

# --- row 18: was CoinExToken -> now One ---
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005700"
$ws.Range("E18").Value = "17OneONE"

# --- more simple price refreshes ---
$ws.Range("D19").Value = "'0.006253"
$ws.Range("D21").Value = "'0.003798"
$ws.Range("D22").Value = "'0.0001497"
$ws.Range("D23").Value = "'0.0004174"
$ws.Range("D24").Value = "'3.969"
$ws.Range("D40").Value = "'0.04188"
$ws.Range("D41").Value = "'0.007004"

# --- row 42: was CEJI -> now BKEXToken ---
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1044"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# --- row 43: was BKEXToken -> now CEJI ---
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002903"
$ws.Range("E43").Value = "42CEJICEJI"

# --- remaining simple price refreshes (rows 44-50) ---
$ws.Range("D44").Value = "'0.009640"
$ws.Range("D45").Value = "'0.00005431"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D47").Value = "'0.6788"
$ws.Range("D48").Value = "'0.03088"
$ws.Range("D49").Value = "'0.00002096"
$ws.Range("D50").Value = "'0.01008"

Write-Host "Applied cryptos.xlsx symbol list update"
